# cryptos.xlsx — refresh the Price (D) and Volume(1h) (E) columns for rows 2-51
# to the latest scrape, matching the GitHub Actions bot's commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "D4" (TetherUSD, fixed at 1.00, never updated) is our reference cell for the
# column's original, un-styled default formatting (no NumberFormat/style override).

# Some new prices read as plain numbers (e.g. "242.83"). Excel would silently
# convert those to numeric cells, but the source data is text (see "36.454.69",
# "1.938.08", etc. elsewhere in the column), so force text with a temporary
# "@" format, assign the value, then restore the cell to the column's normal
# (un-formatted) style so no stray number-format lingers on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.83"
$ws.Range("D5").Style = $ws.Range("D4").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = $ws.Range("D4").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.84"
$ws.Range("D8").Style = $ws.Range("D4").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.364"
$ws.Range("D9").Style = $ws.Range("D4").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.34"
$ws.Range("D10").Style = $ws.Range("D4").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("D11").Style = $ws.Range("D4").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.819"
$ws.Range("D13").Style = $ws.Range("D4").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.35"
$ws.Range("D14").Style = $ws.Range("D4").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.43"
$ws.Range("D16").Style = $ws.Range("D4").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.22"
$ws.Range("D17").Style = $ws.Range("D4").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.31"
$ws.Range("D20").Style = $ws.Range("D4").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.75"
$ws.Range("D22").Style = $ws.Range("D4").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.98"
$ws.Range("D23").Style = $ws.Range("D4").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = $ws.Range("D4").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("D27").Style = $ws.Range("D4").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.23"
$ws.Range("D28").Style = $ws.Range("D4").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.36"
$ws.Range("D29").Style = $ws.Range("D4").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.66"
$ws.Range("D33").Style = $ws.Range("D4").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0627"
$ws.Range("D34").Style = $ws.Range("D4").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.28"
$ws.Range("D35").Style = $ws.Range("D4").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.01"
$ws.Range("D37").Style = $ws.Range("D4").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("D38").Style = $ws.Range("D4").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.14"
$ws.Range("D39").Style = $ws.Range("D4").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = $ws.Range("D4").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.76"
$ws.Range("D45").Style = $ws.Range("D4").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.46"
$ws.Range("D48").Style = $ws.Range("D4").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.17"
$ws.Range("D49").Style = $ws.Range("D4").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.55"
$ws.Range("D51").Style = $ws.Range("D4").Style

# Everything else (prices that already read as text, plus the whole Volume(1h)
# percentage column, which is always padded text like "  -0.39%  ") can be
# assigned directly.
$ws.Range("D2").Value = "36.454.69"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.938.08"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -9.56%  "
$ws.Range("E9").Value = "  -5.72%  "
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -7.42%  "
$ws.Range("E14").Value = "  -7.15%  "
$ws.Range("D15").Value = "2.213.98"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("D18").Value = "1.932.20"
$ws.Range("E18").Value = "  -4.11%  "
$ws.Range("D19").Value = "36.340.16"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  -7.26%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  -7.86%  "
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("E30").Value = "  -8.43%  "
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("E33").Value = "  -7.17%  "
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  -7.17%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  -8.78%  "
$ws.Range("E40").Value = "  -10.14%  "
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  -7.26%  "
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -6.45%  "
$ws.Range("E46").Value = "  -8.40%  "
$ws.Range("D47").Value = "1.343.62"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("E48").Value = "  -8.46%  "
$ws.Range("E49").Value = "  -6.93%  "
$ws.Range("E50").Value = "  -3.00%  "
$ws.Range("E51").Value = "  +1.99%  "
